# Update poll data: insert new SMS Morgan poll rows into the Data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The sheet currently ends at row 141 (MidDate 44142 / Galaxy poll).
# Two new, earlier-dated rows need to be inserted before it, and one new,
# later-dated row needs to be inserted after it, so the data stays sorted
# by MidDate:
#   141 -> 16/09/2020 SMS Morgan   (new)
#   142 -> 30/09/2020 SMS Morgan   (new)
#   143 -> 07/11/2020 Galaxy       (existing row, shifted down from 141)
#   144 -> 10/11/2020 SMS Morgan   (new)

# Push the existing last row (141) down by two rows.
$ws.Rows("141:142").Insert()

# Make room for the new trailing row after the (now) row 143.
$ws.Rows("144:144").Insert()

function Set-PollRow($Row, $MidDate, $Firm, $AlpFp, $LnpFp, $GrnFp, $OthFp, $TppAlp) {
    $ws.Cells.Item($Row, 1).Value = $MidDate
    $ws.Cells.Item($Row, 2).Value = $Firm
    $ws.Cells.Item($Row, 3).Value = $AlpFp
    $ws.Cells.Item($Row, 4).Value = $LnpFp
    $ws.Cells.Item($Row, 5).Value = $GrnFp
    $ws.Cells.Item($Row, 6).Value = $OthFp
    $ws.Cells.Item($Row, 7).Value = "#N/A"
    $ws.Cells.Item($Row, 8).Value = $TppAlp
}

Set-PollRow 141 44090 "SMS Morgan" 51.5 38.5 37 12 12.5
Set-PollRow 142 44104 "SMS Morgan" 51.5 39.5 39 10 11.5
Set-PollRow 144 44145 "SMS Morgan" 58.5 34.5 45 11 9.5

# Match the view state left behind after the edit: the active selection in
# the lower (scrolled) pane ends up on the new blank row below the data.
$ws.Range("A145").Select() | Out-Null

Write-Output "Inserted rows 141, 142 and 144; sheet now spans to row $($ws.UsedRange.Rows.Count)."
